# Commit: "test name modification in input file"
#
# The "MoisturizerTest" test case was renamed to "BuyBodyLotionTest", and the
# "SunscreenTest" test case (and its data block) was removed entirely from
# both the TestStatus and TestCases sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "TestStatus"
$ws2 = $wb.Worksheets.Item(2)   # "TestCases"

# --- TestStatus sheet ---
# Row 3 held "MoisturizerTest" -> rename it; row 4 held "SunscreenTest" -> remove it.
$ws1.Range("A3").Value = "BuyBodyLotionTest"
$ws1.Rows.Item(4).Delete()

# --- TestCases sheet ---
# A5 held the "MoisturizerTest" section header -> rename it.
$ws2.Range("A5").Value = "BuyBodyLotionTest"
# Rows 9-11 held the "SunscreenTest" section (header + data rows) -> remove them.
$ws2.Range("A9:K11").Delete() | Out-Null

# Column A on both sheets needs to widen slightly to fit "BuyBodyLotionTest".
$ws1.Columns.Item(1).ColumnWidth = 15.998697916666666
$ws2.Columns.Item(1).ColumnWidth = 15.998697916666666

# Selection on the TestCases sheet moved to A5.
$ws2.Range("A5").Select() | Out-Null
